$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.631.99'
$ws.Range("E2").Value = '  +4.12%  '

$ws.Range("D3").Value = '1.747.51'
$ws.Range("E3").Value = '  +4.63%  '

$ws.Range("D5").Value = '''247.30'
$ws.Range("E5").Value = '  +3.35%  '

$ws.Range("D7").Value = '''0.4806'
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").Value = '''0.2704'
$ws.Range("E8").Value = '  +2.75%  '

$ws.Range("D9").Value = '''0.06261'
$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").Value = '1.745.84'
$ws.Range("E10").Value = '  +4.51%  '

$ws.Range("E11").Value = '  +1.68%  '

$ws.Range("D12").Value = '''15.83'
$ws.Range("E12").Value = '  +6.43%  '

$ws.Range("D13").Value = '''0.6179'
$ws.Range("E13").Value = '  +4.99%  '

$ws.Range("D14").Value = '''4.510'
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").Value = '''77.32'
$ws.Range("E15").Value = '  +2.66%  '

$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").Value = '26.628.85'
$ws.Range("E17").Value = '  +4.16%  '

$ws.Range("D18").Value = '''1.000'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").Value = '''0.000006909'
$ws.Range("E19").Value = '  +2.20%  '

$ws.Range("D20").Value = '''11.71'
$ws.Range("E20").Value = '  +2.22%  '

$ws.Range("D21").Value = '1.970.71'
$ws.Range("E21").Value = '  +4.58%  '

$ws.Range("D22").Value = '''4.657'
$ws.Range("E22").Value = '  +4.97%  '

$ws.Range("D23").Value = '''8.863'
$ws.Range("E23").Value = '  +1.30%  '

$ws.Range("D24").Value = '''5.357'
$ws.Range("E24").Value = '  +1.51%  '

$ws.Range("D25").Value = '''136.31'
$ws.Range("E25").Value = '  -0.30%  '

$ws.Range("E26").Value = '  +2.93%  '

$ws.Range("D27").Value = '''1.827'
$ws.Range("E27").Value = '  +6.16%  '

$ws.Range("D28").Value = '''1.422'
$ws.Range("E28").Value = '  +2.24%  '

$ws.Range("D29").Value = '''107.59'
$ws.Range("E29").Value = '  +2.62%  '

$ws.Range("D30").Value = '''4.029'
$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").Value = '''3.771'
$ws.Range("E31").Value = '  +3.37%  '

$ws.Range("D32").Value = '''0.07901'
$ws.Range("E32").Value = '  +0.92%  '

$ws.Range("D33").Value = '''0.04580'
$ws.Range("E33").Value = '  +8.40%  '

$ws.Range("D34").Value = '''2.614'
$ws.Range("E34").Value = '  -0.20%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.001'
$ws.Range("E35").Value = '  +4.95%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.6356'
$ws.Range("E36").Value = '  +4.27%  '

$ws.Range("D37").Value = '''0.9576'
$ws.Range("E37").Value = '  +11.32%  '

$ws.Range("D38").Value = '''114.73'
$ws.Range("E38").Value = '  +18.86%  '

$ws.Range("D39").Value = '''2.478'
$ws.Range("E39").Value = '  -4.52%  '

$ws.Range("D40").Value = '''1.979'
$ws.Range("E40").Value = '  +5.52%  '

$ws.Range("D41").Value = '''1.004'
$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("E42").Value = '  +2.92%  '

$ws.Range("D43").Value = '''5.713'
$ws.Range("E43").Value = '  +17.33%  '

$ws.Range("D44").Value = '''0.3922'
$ws.Range("E44").Value = '  +4.12%  '

$ws.Range("D45").Value = '''6.765'
$ws.Range("E45").Value = '  +8.71%  '

$ws.Range("D46").Value = '''0.1206'
$ws.Range("E46").Value = '  +7.92%  '

$ws.Range("D47").Value = '''0.05327'
$ws.Range("E47").Value = '  +1.22%  '

$ws.Range("D48").Value = '''7.965'
$ws.Range("E48").Value = '  +7.92%  '

$ws.Range("D49").Value = '''30.86'
$ws.Range("E49").Value = '  +3.22%  '

$ws.Range("E50").Value = '  +3.52%  '

$ws.Range("D51").Value = '''51.77'
$ws.Range("E51").Value = '  +3.49%  '
